{"js": "// The document is a table of \"three-digit \u00f7 one-digit\" division practice\n// problems (e.g. \"265\u00f76=44, 1\"). Each cell's full text is replaced with a\n// freshly generated problem string. Every \"old\" value below is unique and\n// appears exactly once in the document, so body.search() for the exact,\n// full old text unambiguously finds the single run/paragraph to update.\nconst replacements = [\n  [\"265\u00f76=44, 1\", \"632\u00f74=158, 0\"],\n  [\"796\u00f73=265, 1\", \"324\u00f79=36, 0\"],\n  [\"660\u00f75=132, 0\", \"175\u00f77=25, 0\"],\n  [\"828\u00f75=165, 3\", \"344\u00f78=43, 0\"],\n  [\"469\u00f73=156, 1\", \"175\u00f76=29, 1\"],\n  [\"465\u00f72=232, 1\", \"853\u00f78=106, 5\"],\n  [\"498\u00f77=71, 1\", \"233\u00f75=46, 3\"],\n  [\"920\u00f75=184, 0\", \"451\u00f79=50, 1\"],\n  [\"705\u00f79=78, 3\", \"422\u00f77=60, 2\"],\n  [\"642\u00f77=91, 5\", \"668\u00f74=167, 0\"],\n  [\"146\u00f72=73, 0\", \"865\u00f73=288, 1\"],\n  [\"602\u00f78=75, 2\", \"697\u00f77=99, 4\"],\n  [\"687\u00f75=137, 2\", \"494\u00f75=98, 4\"],\n  [\"616\u00f73=205, 1\", \"231\u00f79=25, 6\"],\n  [\"731\u00f76=121, 5\", \"624\u00f77=89, 1\"],\n  [\"147\u00f73=49, 0\", \"220\u00f73=73, 1\"],\n  [\"113\u00f78=14, 1\", \"775\u00f73=258, 1\"],\n  [\"948\u00f76=158, 0\", \"120\u00f78=15, 0\"],\n  [\"857\u00f77=122, 3\", \"100\u00f73=33, 1\"],\n  [\"119\u00f75=23, 4\", \"974\u00f78=121, 6\"],\n  [\"532\u00f77=76, 0\", \"152\u00f72=76, 0\"],\n  [\"841\u00f75=168, 1\", \"566\u00f75=113, 1\"],\n  [\"931\u00f76=155, 1\", \"826\u00f72=413, 0\"],\n  [\"140\u00f76=23, 2\", \"248\u00f72=124, 0\"],\n  [\"153\u00f73=51, 0\", \"620\u00f78=77, 4\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each three-digit division equation with its new value using\n# Word's Find/Replace (Range.Find), matching the exact, unique cell text.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"265\u00f76=44, 1\", \"632\u00f74=158, 0\"),\n  @(\"796\u00f73=265, 1\", \"324\u00f79=36, 0\"),\n  @(\"660\u00f75=132, 0\", \"175\u00f77=25, 0\"),\n  @(\"828\u00f75=165, 3\", \"344\u00f78=43, 0\"),\n  @(\"469\u00f73=156, 1\", \"175\u00f76=29, 1\"),\n  @(\"465\u00f72=232, 1\", \"853\u00f78=106, 5\"),\n  @(\"498\u00f77=71, 1\", \"233\u00f75=46, 3\"),\n  @(\"920\u00f75=184, 0\", \"451\u00f79=50, 1\"),\n  @(\"705\u00f79=78, 3\", \"422\u00f77=60, 2\"),\n  @(\"642\u00f77=91, 5\", \"668\u00f74=167, 0\"),\n  @(\"146\u00f72=73, 0\", \"865\u00f73=288, 1\"),\n  @(\"602\u00f78=75, 2\", \"697\u00f77=99, 4\"),\n  @(\"687\u00f75=137, 2\", \"494\u00f75=98, 4\"),\n  @(\"616\u00f73=205, 1\", \"231\u00f79=25, 6\"),\n  @(\"731\u00f76=121, 5\", \"624\u00f77=89, 1\"),\n  @(\"147\u00f73=49, 0\", \"220\u00f73=73, 1\"),\n  @(\"113\u00f78=14, 1\", \"775\u00f73=258, 1\"),\n  @(\"948\u00f76=158, 0\", \"120\u00f78=15, 0\"),\n  @(\"857\u00f77=122, 3\", \"100\u00f73=33, 1\"),\n  @(\"119\u00f75=23, 4\", \"974\u00f78=121, 6\"),\n  @(\"532\u00f77=76, 0\", \"152\u00f72=76, 0\"),\n  @(\"841\u00f75=168, 1\", \"566\u00f75=113, 1\"),\n  @(\"931\u00f76=155, 1\", \"826\u00f72=413, 0\"),\n  @(\"140\u00f76=23, 2\", \"248\u00f72=124, 0\"),\n  @(\"153\u00f73=51, 0\", \"620\u00f78=77, 4\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n  if (-not $found) {\n    throw \"No match found for: $oldText\"\n  }\n}\n"}
